# Populate the "funcionarios" (employees) sheet with the roster data,
# apply the date/header number formatting, size the columns, mark the
# last-used cell (H8) with the lingering underline format, and leave the
# selection where the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Employee roster -------------------------------------------------
# code, name, admission date (as an Excel 1900-date-system serial number,
# so it round-trips as a clean integer with no time-of-day component),
# setor, funcao
$data = @(
    @(1000, "João Silva",        44941, "Recursos Humanos",              "Analista de Recursos Humanos"),
    @(1001, "Maria Santos",      44977, "Financeiro",                    "Analista Financeiro"),
    @(1002, "Pedro Oliveira",    44995, "Marketing",                     "Gerente de Marketing"),
    @(1003, "Ana Souza",         45021, "Vendas",                        "Executivo de Vendas"),
    @(1004, "Lucas Pereira",     45058, "Produção",                      "Supervisor de Produção"),
    @(1005, "Camila Almeida",    45102, "TI (Tecnologia da Informação)", "Desenvolvedor de Software"),
    @(1006, "Felipe Rodrigues",  45125, "Logística",                     "Coordenador de Logística"),
    @(1007, "Juliana Lima",      45147, "Qualidade",                     "Especialista em Controle de Qualidade"),
    @(1008, "Rafael Costa",      45199, "Jurídico",                      "Advogado"),
    @(1009, "Larissa Fernandes", 45220, "Administrativo",                "Assistente Administrativo")
)

$row = 2
foreach ($rec in $data) {
    $code   = $rec[0]
    $name   = $rec[1]
    $serial = $rec[2]
    $setor  = $rec[3]
    $funcao = $rec[4]

    $ws.Cells.Item($row, 1).Value = $code
    $ws.Cells.Item($row, 2).Value = $name
    $ws.Cells.Item($row, 3).Value = $serial
    $ws.Cells.Item($row, 4).Value = $setor
    $ws.Cells.Item($row, 5).Value = $funcao

    $row++
}

# --- Number formatting -------------------------------------------------
# data_admissao column (C) and its header both use the short-date format.
$ws.Range("C2:C11").NumberFormat = "mm-dd-yy"
$ws.Range("C1").NumberFormat = "mm-dd-yy"

# --- Column widths -------------------------------------------------
$ws.Columns("A").ColumnWidth = 23.14
$ws.Columns("B").ColumnWidth = 22.43
$ws.Columns("C").ColumnWidth = 19.71
$ws.Columns("D").ColumnWidth = 27.57
$ws.Columns("E").ColumnWidth = 19.86

# --- Stray formatted cell left by the author (H8) -------------------------------------------------
$ws.Cells.Item(8, 8).Font.Underline = $true

# --- Selection, matching where the author's cursor ended up -------------------------------------------------
$ws.Range("H8").Select() | Out-Null

Write-Output "employee roster populated"
